$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the last three "Date" column values with their literal text representation
$ws.Range("C7").Value = "17/06/1997"
$ws.Range("C8").Value = "18/06/1997"
$ws.Range("C9").Value = "19/06/1997"

# Add a new row of data (times past midnight / a later date)
$ws.Range("A10").Value = 1.61111111111111
$ws.Range("B10").Value = 1.62569444444444
$ws.Range("C10").Value = 50575
$ws.Range("A10:B10").NumberFormat = "hh:mm:ss"
$ws.Range("C10").NumberFormat = "dd/mm/yy"

# Widen the Date column so the new values are fully visible
$ws.Columns.Item(3).ColumnWidth = 28.43

$null = $ws.Range("F12").Select()
